{"js": "// Esempio didattico quasi completo\n//\n// The three bullet paragraphs rotate their content:\n//   \"Provare a fare un sistemino...\"                 -> \"Perche' non posso togliere...\"\n//   \"Perche' non posso togliere...\" (+ bookmark)      -> \"Errori di estensione con il doppio if: controllare\"\n//   \"Svecchiare anche il doppio blocco if\"            -> \"Buildare il\" + _GoBack bookmark + \" sistemino...\"\n//\n// Locate the paragraphs by their current text (robust to ordering/index\n// drift) instead of hard-coded paragraph indices.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet pProvare = null;\nlet pPerche = null;\nlet pSvecchiare = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t.indexOf(\"Provare a fare un sistemino\") !== -1) {\n    pProvare = paras.items[i];\n  } else if (t.indexOf(\"Perche\\u2019 non posso togliere\") !== -1) {\n    pPerche = paras.items[i];\n  } else if (t.indexOf(\"Svecchiare anche il doppio blocco if\") !== -1) {\n    pSvecchiare = paras.items[i];\n  }\n}\n\nif (!pProvare || !pPerche || !pSvecchiare) {\n  throw new Error(\"Could not locate the expected source paragraphs.\");\n}\n\n// 1) \"Provare a fare...\" paragraph becomes the old \"Perche' non posso...\" text.\npProvare.insertText(\n  \"Perche\\u2019 non posso togliere un pezzo da sotto un blocco e attaccarlo altrove?\",\n  \"Replace\"\n);\n\n// 2) \"Perche' non posso...\" paragraph becomes the new \"Errori di estensione...\" text.\npPerche.insertText(\n  \"Errori di estensione con il doppio if: controllare\",\n  \"Replace\"\n);\n\n// 3) \"Svecchiare...\" paragraph becomes \"Buildare il\" (first half for now).\npSvecchiare.insertText(\"Buildare il\", \"Replace\");\n\nawait context.sync();\n\n// The _GoBack bookmark used to sit at the end of the \"Perche' non posso...\"\n// paragraph; remove it from there before recreating it in its new home.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Append the second half of the text onto the (now plain) paragraph.\npSvecchiare\n  .getRange(\"End\")\n  .insertText(\n    \" sistemino da mandare a Denti che permetta di interpretare semplici frasi.\",\n    \"End\"\n  );\n\nawait context.sync();\n\n// Re-locate the \"Buildare il\" run via search so we get a freshly anchored,\n// non-stale range, then drop the _GoBack bookmark right after it \u2014 between\n// the two runs, exactly where it lived (relative to the split) originally.\nconst searchResults = pSvecchiare.search(\"Buildare il\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nsearchResults.items[0].getRange(\"End\").insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Esempio didattico quasi completo\n#\n# The three bullet paragraphs rotate their content:\n#   \"Provare a fare un sistemino...\"                 -> \"Perche' non posso togliere...\"\n#   \"Perche' non posso togliere...\" (+ bookmark)      -> \"Errori di estensione con il doppio if: controllare\"\n#   \"Svecchiare anche il doppio blocco if\"            -> \"Buildare il\" + _GoBack bookmark + \" sistemino...\"\n\n$d = $word.ActiveDocument\n$apos = [char]0x2019\n\n# Locate the three target paragraphs by their current text (robust to\n# index/ordering drift) instead of relying on hard-coded paragraph numbers.\n$idxProvare = $null\n$idxPerche = $null\n$idxSvecchiare = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"Provare a fare un sistemino\")) { $idxProvare = $i }\n    elseif ($t.Contains(\"Perche\" + $apos + \" non posso togliere\")) { $idxPerche = $i }\n    elseif ($t.Contains(\"Svecchiare anche il doppio blocco if\")) { $idxSvecchiare = $i }\n}\n\nif (-not $idxProvare -or -not $idxPerche -or -not $idxSvecchiare) {\n    throw \"Could not locate the expected source paragraphs.\"\n}\n\n# 1) \"Provare a fare...\" paragraph becomes the old \"Perche' non posso...\" text.\n$d.Paragraphs.Item($idxProvare).Range.Text = \"Perche\" + $apos + \" non posso togliere un pezzo da sotto un blocco e attaccarlo altrove?\"\n\n# 2) \"Perche' non posso...\" paragraph becomes the new \"Errori di estensione...\" text.\n$d.Paragraphs.Item($idxPerche).Range.Text = \"Errori di estensione con il doppio if: controllare\"\n\n# The _GoBack bookmark used to sit at the end of that paragraph's text;\n# remove it from there before recreating it in its new home below.\n$d.Bookmarks.Item(\"_GoBack\").Delete() | Out-Null\n\n# 3) \"Svecchiare...\" paragraph becomes \"Buildare il\" + bookmark + rest of text.\n$d.Paragraphs.Item($idxSvecchiare).Range.Text = \"Buildare il\"\n$d.Paragraphs.Item($idxSvecchiare).Range.InsertAfter(\" sistemino da mandare a Denti che permetta di interpretare semplici frasi.\")\n\n# Re-find the \"Buildare il\" run so we get a freshly anchored range (avoids\n# stale collapsed-range/paragraph-boundary quirks), then drop the _GoBack\n# bookmark right after it -- between the two runs.\n$searchRange = $d.Paragraphs.Item($idxSvecchiare).Range.Duplicate\n$searchRange.Find.Execute(\"Buildare il\") | Out-Null\n$bmPoint = $d.Range($searchRange.End, $searchRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmPoint) | Out-Null\n"}
